$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 315, shifting existing rows 315:342 down to 316:343
$ws.Rows.Item(315).Insert()

# Populate the new row 315 with the new weekly data point (copy of the
# constant columns plus the new date/volume/price figures)
$ws.Cells.Item(315, 1).Value = 6
$ws.Cells.Item(315, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(315, 3).Value = "Metropolitana"
$ws.Cells.Item(315, 4).Value = 45106
$ws.Cells.Item(315, 4).Style = $ws.Cells.Item(316, 4).Style
$ws.Cells.Item(315, 4).NumberFormat = $ws.Cells.Item(316, 4).NumberFormat
$ws.Cells.Item(315, 5).Value = 13
$ws.Cells.Item(315, 6).Value = 100112029
$ws.Cells.Item(315, 7).Value = "Orégano"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 38
$ws.Cells.Item(315, 11).Value = 19000
$ws.Cells.Item(315, 12).Value = 20000
$ws.Cells.Item(315, 13).Value = 19526
$ws.Cells.Item(315, 14).Value = "$/docena de atados"
$ws.Cells.Item(315, 15).Value = "Región Metropolitana"
$ws.Cells.Item(315, 16).Value = 6509
$ws.Cells.Item(315, 17).Value = 3
$ws.Cells.Item(315, 18).Value = "Hortaliza"
